# Applies the commit's changes to interested_tickers_oanda.xlsx:
#  - sheet "h4" gets the full ticker list (same order as sheets "d"/"m1"/"h1")
#    added into A2:A10 (it previously only had EUR_USD/GBP_USD in A2:A3).
#  - sheet "h4"'s duplicate-values conditional formatting is rebuilt to
#    cover A2:A10 instead of A4:A10 (this creates a new dxf entry, which is
#    why every sheet's conditional-format dxfId shifts up by one).
#  - sheet "d" becomes the active/selected tab (previously "h4" was), with
#    its last selected cell moving from D21 to D24.
#  - sheet "h4"'s own last selected cell becomes C2 (no longer the active tab).
#  - sheet "h1"'s last selection changes from the whole A1:A10 to A2:A10.

$wb = $excel.ActiveWorkbook

$wsD  = $wb.Worksheets.Item("d")
$wsM1 = $wb.Worksheets.Item("m1")
$wsH1 = $wb.Worksheets.Item("h1")
$wsH4 = $wb.Worksheets.Item("h4")

# --- sheet "h4": fill in the rest of the ticker list -----------------------
$wsH4.Range("A2").Value  = "AUD_USD"
$wsH4.Range("A3").Value  = "EUR_JPY"
$wsH4.Range("A4").Value  = "EUR_USD"
$wsH4.Range("A5").Value  = "GBP_JPY"
$wsH4.Range("A6").Value  = "GBP_USD"
$wsH4.Range("A7").Value  = "NZD_USD"
$wsH4.Range("A8").Value  = "USD_CAD"
$wsH4.Range("A9").Value  = "USD_CHF"
$wsH4.Range("A10").Value = "USD_JPY"

# --- sheet "h4": move conditional formatting from A4:A10 to A2:A10 ---------
$wsH4.Range("A4:A10").FormatConditions.Delete()
$dupRule = $wsH4.Range("A2:A10").FormatConditions.AddUniqueValues()
$dupRule.DupeUnique = 1
$dupRule.Interior.Color = 13551615
$dupRule.Font.Color = 393372

# --- per-sheet selections ---------------------------------------------------
$wsH4.Activate() | Out-Null
$wsH4.Range("C2").Select() | Out-Null

$wsH1.Activate() | Out-Null
$wsH1.Range("A2:A10").Select() | Out-Null

# "d" ends up the active/selected tab (was "h4" before the edit)
$wsD.Activate() | Out-Null
$wsD.Range("D24").Select() | Out-Null

Write-Host "Applied interested_tickers_oanda.xlsx update"
